$d = $word.ActiveDocument

function Merge-Paragraph {
    param($doc, $index, $newText)
    $target = $doc.Paragraphs.Item($index)
    $target.Range.InsertParagraphBefore()
    $newPara = $doc.Paragraphs.Item($index)
    $newRange = $doc.Range($newPara.Range.Start, $newPara.Range.End - 1)
    $newRange.Text = $newText
    $oldPara = $doc.Paragraphs.Item($index + 1)
    $oldFull = $doc.Range($oldPara.Range.Start, $oldPara.Range.End)
    $oldFull.Delete()
}

# 1) Cat/Parrot/Seed solution paragraph: merge the two runs that were split by a
#    gramStart/gramEnd proofErr pair into a single run.
Merge-Paragraph $d 15 "The solution would keep the parrot and cat separate and also keep the seeds and parrot separate.  The first move would be for the man to transport the parrot across the river.  This satisfies both requirements.  The man would travel back alone and pick up either the cat or the seeds.  He would then transport the item across and pick up the parrot because the parrot cannot be left alone on either bank.  He travels across the river with the parrot to drop it off and pick up the item left on the bank.  He travels back across and drops off the item.  He then travels back alone to pick up the parrot and finally transport it across the river."

# 2) Socks question paragraph: same kind of merge.
Merge-Paragraph $d 18 "There are 20 socks in a drawer: 5 pairs of black socks, 3 pairs of brown and 2 pairs of white.  You select the socks in the dark and can check them only after a selection has been made.  What is the smallest number of socks you need to select to guarantee getting the following:"

# 3) / 4) / 5) the girl-counting a)/b)/c) paragraphs.
Merge-Paragraph $d 36 "a) What if the girl counts from 1 to 10"
Merge-Paragraph $d 37 "b) What if the girl counts from 1 to 100"
Merge-Paragraph $d 38 "c) What if the girl counts from 1 to 1000"

# 6) Answer the third question's "What are potential solutions?" part: insert a new
#    paragraph with the answer text, and move the _GoBack bookmark so it ends up at
#    the end of that new paragraph (where the last edit now is), replacing the old
#    trailing empty paragraph.
$p43 = $d.Paragraphs.Item(43)
$p45 = $d.Paragraphs.Item(45)
$range = $d.Range($p43.Range.Start, $p45.Range.End)

$frag = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:r><w:t>A challenge this problem faces is that the solution should work for any number.</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>What are potential solutions?</w:t></w:r></w:p>
<w:p><w:r><w:t>A potential solution is to separate the initial rotation, which is 9, from the rest, which is 8.  You could also break down her counting into a simpler pattern of 5 to start then 4.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$range.InsertXML($frag) | Out-Null
